# resolve domain related templates
#
# The "Domains" worksheet had a "Length" column (column F) that is no
# longer needed - remove it entirely (cells shift left, the "Length"
# shared string is dropped) and update the active selection to reflect
# where the user ended up after the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Domains")

# Delete the whole "Length" column (F); everything to its right (Nullable,
# Unique, Default Value, Comment) shifts one column to the left.
$ws.Columns("F").Delete()

# Deleting the column nudges the frozen-pane selections; restore them to
# match what was recorded after the edit (top-right pane stays at B1, the
# active/bottom-right pane moves to H8).
$ws.Range("B1").Select() | Out-Null
$ws.Range("H8").Select() | Out-Null
